$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the phone number on row 2 (drop stray leading "55")
$ws.Range("B2").Value = 41998306017

# Remove the extra rows (Victor Cals / Guilherme) - no longer needed
$ws.Range("A3:B4").EntireRow.Delete() | Out-Null

# Narrow column B back down now that the long rows are gone
$ws.Columns.Item(2).ColumnWidth = 11.5

# Leave the selection where the last edit happened
$ws.Range("B3").Select() | Out-Null
